$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "value" column (column C) from 0.05 to 0.04 for rows 2 through 402
$ws.Range("C2:C402").Value = 0.04

# Update the active cell selection to C9 (was C10)
$ws.Range("C9").Select()
